$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2458.5
$ws.Range("J40").Value = 3316.6667
$ws.Range("L40").Value = 3316.6667
$ws.Range("N40").Value = -3666.6667

$ws.Range("H48").Value = 6500
$ws.Range("I48").Value = 6583.3335
$ws.Range("K48").Value = 19750.0005
$ws.Range("M48").Value = -19458.0005

$ws.Range("H56").Value = 6500
$ws.Range("I56").Value = 6583.3335
$ws.Range("K56").Value = 19750.0005
$ws.Range("M56").Value = -19216.0005

$ws.Range("H64").Value = 7107.143
$ws.Range("I64").Value = 6916.6665
$ws.Range("J64").Value = 7250
$ws.Range("K64").Value = 6916.6665
$ws.Range("L64").Value = 7250
$ws.Range("M64").Value = -6668.6665
$ws.Range("N64").Value = -7746

$ws.Range("H67").Value = 7107.143
$ws.Range("I67").Value = 6916.6665
$ws.Range("J67").Value = 7250
$ws.Range("K67").Value = 6916.6665
$ws.Range("L67").Value = 7250
$ws.Range("M67").Value = -6058.6665
$ws.Range("N67").Value = -8966

$ws.Range("H107").Value = 1080.8529
$ws.Range("I107").Value = 935.7778
$ws.Range("J107").Value = 1640.4286
$ws.Range("K107").Value = 935.7778
$ws.Range("L107").Value = 1640.4286
$ws.Range("M107").Value = 984.2222
$ws.Range("N107").Value = -5480.4286

$ws.Range("H116").Value = 18889.37
$ws.Range("I116").Value = 23921
$ws.Range("K116").Value = 23921
$ws.Range("M116").Value = -20479

$ws.Range("H132").Value = 28984.076
$ws.Range("I132").Value = 31011.916
$ws.Range("K132").Value = 93035.74800000001
$ws.Range("M132").Value = -90505.74800000001

$ws.Range("H135").Value = 2178.5
$ws.Range("I135").Value = 1772.5
$ws.Range("J135").Value = 3498
$ws.Range("K135").Value = 15952.5
$ws.Range("L135").Value = 31482
$ws.Range("M135").Value = -13417.5
$ws.Range("N135").Value = -36552

$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

$ws.Range("H141").Value = 1908.6666
$ws.Range("I141").Value = 1908.6666
$ws.Range("K141").Value = 5725.9998
$ws.Range("M141").Value = -545.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200.33333
$ws.Range("I4").Value = 250.5
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 250.5
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -134.5
$ws.Range("N4").Value = -332

$ws.Range("H110").Value = 4478.9214
$ws.Range("I110").Value = 4735.364
$ws.Range("K110").Value = 4735.364
$ws.Range("M110").Value = -2690.364

$ws.Range("H122").Value = 2081.2942
$ws.Range("I122").Value = 1861.3572
$ws.Range("J122").Value = 3107.6667
$ws.Range("K122").Value = 5584.071599999999
$ws.Range("L122").Value = 9323.000100000001
$ws.Range("M122").Value = -3134.071599999999
$ws.Range("N122").Value = -14223.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1271
$ws.Range("J64").Value = 1166.25
$ws.Range("L64").Value = 1166.25
$ws.Range("N64").Value = -1616.25

$ws.Range("H67").Value = 1271
$ws.Range("J67").Value = 1166.25
$ws.Range("L67").Value = 1166.25
$ws.Range("N67").Value = -2726.25

$ws.Range("H105").Value = 4067.6155
$ws.Range("I105").Value = 3870.5557
$ws.Range("K105").Value = 3870.5557
$ws.Range("M105").Value = -2123.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2455.92
$ws.Range("I16").Value = 2317.1052
$ws.Range("K16").Value = 2317.1052
$ws.Range("M16").Value = -2030.1052

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H63").Value = 33750

$ws.Range("H66").Value = 33750

$ws.Range("H94").Value = 2332.389
$ws.Range("I94").Value = 1386
$ws.Range("J94").Value = 3819.5715
$ws.Range("K94").Value = 1386
$ws.Range("L94").Value = 3819.5715
$ws.Range("M94").Value = -935
$ws.Range("N94").Value = -4721.5715

$ws.Range("H105").Value = 45540.668
$ws.Range("I105").Value = 65811.164
$ws.Range("K105").Value = 65811.164
$ws.Range("M105").Value = -64064.164

$ws.Range("H113").Value = 2455.92
$ws.Range("I113").Value = 2317.1052
$ws.Range("K113").Value = 2317.1052
$ws.Range("M113").Value = -147.1052

$ws.Range("H122").Value = 2012.5555
$ws.Range("I122").Value = 1847.5
$ws.Range("K122").Value = 5542.5
$ws.Range("M122").Value = -3092.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 556.46155
$ws.Range("I5").Value = 549.1818
$ws.Range("J5").Value = 596.5
$ws.Range("K5").Value = 1647.5454
$ws.Range("L5").Value = 1789.5
$ws.Range("M5").Value = -1535.5454
$ws.Range("N5").Value = -2013.5

$ws.Range("H135").Value = 556.46155
$ws.Range("I135").Value = 549.1818
$ws.Range("J135").Value = 596.5
$ws.Range("K135").Value = 4942.6362
$ws.Range("L135").Value = 5368.5
$ws.Range("M135").Value = -2407.6362
$ws.Range("N135").Value = -10438.5

$ws.Range("H140").Value = 1879.5
$ws.Range("I140").Value = 1858.2778
$ws.Range("K140").Value = 5574.8334
$ws.Range("M140").Value = -394.8334000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 18430.033
$ws.Range("I102").Value = 25449.096
$ws.Range("J102").Value = 2052.2222
$ws.Range("K102").Value = 25449.096
$ws.Range("L102").Value = 2052.2222
$ws.Range("M102").Value = -23827.096
$ws.Range("N102").Value = -5296.2222

$ws.Range("H113").Value = 2666.88
$ws.Range("I113").Value = 2461.6875
$ws.Range("K113").Value = 2461.6875
$ws.Range("M113").Value = -291.6875

$ws.Range("H122").Value = 4468.0586
$ws.Range("I122").Value = 4242.25
$ws.Range("K122").Value = 12726.75
$ws.Range("M122").Value = -10276.75

$ws.Range("H126").Value = 2437
$ws.Range("I126").Value = 1928
$ws.Range("K126").Value = 5784
$ws.Range("M126").Value = -3314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2022.1471
$ws.Range("I46").Value = 1292.0769
$ws.Range("J46").Value = 2474.0952
$ws.Range("K46").Value = 1292.0769
$ws.Range("L46").Value = 2474.0952
$ws.Range("M46").Value = -1104.0769
$ws.Range("N46").Value = -2850.0952

$ws.Range("H136").Value = 4276.625
$ws.Range("I136").Value = 4707
$ws.Range("J136").Value = 4133.1665
$ws.Range("K136").Value = 14121
$ws.Range("L136").Value = 12399.4995
$ws.Range("M136").Value = -11571
$ws.Range("N136").Value = -17499.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10459.4
$ws.Range("J41").Value = 8349.25
$ws.Range("L41").Value = 8349.25
$ws.Range("N41").Value = -9129.25

$ws.Range("H122").Value = 8954554
$ws.Range("I122").Value = 11142861
$ws.Range("K122").Value = 33428583
$ws.Range("M122").Value = -33426133
